$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 720
$ws.Range("I2").Value = 100
$ws.Range("K2").Value = 100
$ws.Range("M2").Value = 13

# Row 7
$ws.Range("H7").Value = 8000
$ws.Range("I7").Value = 8000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 8000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -7888
$ws.Range("N7").Value = ""

# Row 14
$ws.Range("H14").Value = 8000
$ws.Range("I14").Value = 8000
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 8000
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -7809
$ws.Range("N14").Value = ""

# Row 69
$ws.Range("H69").Value = 7388.0513
$ws.Range("J69").Value = 7388.0513
$ws.Range("L69").Value = 22164.1539
$ws.Range("N69").Value = -23912.1539

# Row 72
$ws.Range("H72").Value = 7388.0513
$ws.Range("J72").Value = 7388.0513
$ws.Range("L72").Value = 66492.4617
$ws.Range("N72").Value = -75228.4617

# Row 106
$ws.Range("H106").Value = 2173.4
$ws.Range("I106").Value = 2173.4
$ws.Range("K106").Value = 2173.4
$ws.Range("M106").Value = -1542.4

# Row 107
$ws.Range("H107").Value = 383.33334
$ws.Range("I107").Value = 232.42857
$ws.Range("K107").Value = 232.42857
$ws.Range("M107").Value = 1687.57143


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3678.963
$ws.Range("I32").Value = 710.087
$ws.Range("K32").Value = 710.087
$ws.Range("M32").Value = -423.087

# Row 93
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = ""

# Row 94
$ws.Range("H94").Value = 65000
$ws.Range("J94").Value = 65000
$ws.Range("L94").Value = 65000
$ws.Range("N94").Value = -66802

# Row 95
$ws.Range("H95").Value = 7885.143
$ws.Range("J95").Value = 7885.143
$ws.Range("L95").Value = 7885.143
$ws.Range("N95").Value = -13377.143

# Row 97
$ws.Range("H97").Value = 1110.75
$ws.Range("I97").Value = 1547
$ws.Range("J97").Value = 674.5
$ws.Range("K97").Value = 1547
$ws.Range("L97").Value = 674.5
$ws.Range("M97").Value = -1051
$ws.Range("N97").Value = -1666.5

# Row 132
$ws.Range("H132").Value = 2071.8667
$ws.Range("I132").Value = 1934.1428
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 5802.428400000001
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -3272.428400000001
$ws.Range("N132").Value = -17060


$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1884.2354
$ws.Range("I86").Value = 1925.6923
$ws.Range("J86").Value = 1749.5
$ws.Range("K86").Value = 1925.6923
$ws.Range("L86").Value = 1749.5
$ws.Range("M86").Value = -802.6922999999999
$ws.Range("N86").Value = -3995.5

# Row 89
$ws.Range("H89").Value = 1884.2354
$ws.Range("I89").Value = 1925.6923
$ws.Range("J89").Value = 1749.5
$ws.Range("K89").Value = 9628.461499999999
$ws.Range("L89").Value = 8747.5
$ws.Range("M89").Value = -4012.461499999999
$ws.Range("N89").Value = -19979.5

# Row 94
$ws.Range("H94").Value = 545.2
$ws.Range("I94").Value = 506.75
$ws.Range("K94").Value = 506.75
$ws.Range("M94").Value = -55.75

# Row 107
$ws.Range("H107").Value = 3041.4167
$ws.Range("I107").Value = 2856.8572
$ws.Range("J107").Value = 3299.8
$ws.Range("K107").Value = 2856.8572
$ws.Range("L107").Value = 3299.8
$ws.Range("M107").Value = -936.8571999999999
$ws.Range("N107").Value = -7139.8


$ws = $wb.Worksheets.Item("CRP")
# Row 10
$ws.Range("H10").Value = 257.125
$ws.Range("I10").Value = 261.4
$ws.Range("J10").Value = 250
$ws.Range("K10").Value = 261.4
$ws.Range("L10").Value = 250
$ws.Range("M10").Value = -122.4
$ws.Range("N10").Value = -528

# Row 16
$ws.Range("H16").Value = 1266.3334
$ws.Range("I16").Value = 1266.3334
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1266.3334
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -979.3334
$ws.Range("N16").Value = ""

# Row 113
$ws.Range("H113").Value = 1266.3334
$ws.Range("I113").Value = 1266.3334
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1266.3334
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 903.6666
$ws.Range("N113").Value = ""


$ws = $wb.Worksheets.Item("CUL")
# Row 49
$ws.Range("H49").Value = 2667.6667
$ws.Range("I49").Value = 3
$ws.Range("J49").Value = 4000
$ws.Range("K49").Value = 9
$ws.Range("L49").Value = 12000
$ws.Range("M49").Value = 147
$ws.Range("N49").Value = -12312

# Row 55
$ws.Range("H55").Value = 4796.6665
$ws.Range("J55").Value = 4796.6665
$ws.Range("L55").Value = 14389.9995
$ws.Range("N55").Value = -14743.9995

# Row 68
$ws.Range("H68").Value = 567.6667
$ws.Range("I68").Value = 300
$ws.Range("J68").Value = 701.5
$ws.Range("K68").Value = 900
$ws.Range("L68").Value = 2104.5
$ws.Range("M68").Value = -89
$ws.Range("N68").Value = -3726.5

# Row 71
$ws.Range("H71").Value = 567.6667
$ws.Range("I71").Value = 300
$ws.Range("J71").Value = 701.5
$ws.Range("K71").Value = 2700
$ws.Range("L71").Value = 6313.5
$ws.Range("M71").Value = 1356
$ws.Range("N71").Value = -14425.5

# Row 80
$ws.Range("H80").Value = 4282.36
$ws.Range("I80").Value = 4012.5264
$ws.Range("K80").Value = 12037.5792
$ws.Range("M80").Value = -11101.5792

# Row 83
$ws.Range("H83").Value = 4282.36
$ws.Range("I83").Value = 4012.5264
$ws.Range("K83").Value = 36112.7376
$ws.Range("M83").Value = -31432.7376

# Row 103
$ws.Range("H103").Value = 1630.1538
$ws.Range("J103").Value = 1630.1538
$ws.Range("L103").Value = 4890.4614
$ws.Range("N103").Value = -6648.4614

# Row 137
$ws.Range("H137").Value = 4028.5715
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 4028.5715
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 12085.7145
$ws.Range("M137").Value = ""
$ws.Range("N137").Value = -22285.7145

# Row 140
$ws.Range("H140").Value = 3039.9167
$ws.Range("I140").Value = 2587.9
$ws.Range("K140").Value = 7763.700000000001
$ws.Range("M140").Value = -2583.700000000001


$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 1920.1
$ws.Range("I122").Value = 1920.1
$ws.Range("K122").Value = 5760.299999999999
$ws.Range("M122").Value = -3310.299999999999

# Row 126
$ws.Range("H126").Value = 3110.889
$ws.Range("I126").Value = 3110.889
$ws.Range("K126").Value = 9332.667000000001
$ws.Range("M126").Value = -6862.667000000001


$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1619
$ws.Range("J22").Value = 1922
$ws.Range("L22").Value = 1922
$ws.Range("N22").Value = -2512

# Row 27
$ws.Range("H27").Value = 1619
$ws.Range("J27").Value = 1922
$ws.Range("L27").Value = 1922
$ws.Range("N27").Value = -2136

# Row 40
$ws.Range("H40").Value = 2776.7778
$ws.Range("I40").Value = 2640.4
$ws.Range("J40").Value = 2947.25
$ws.Range("K40").Value = 2640.4
$ws.Range("L40").Value = 2947.25
$ws.Range("M40").Value = -2504.4
$ws.Range("N40").Value = -3219.25

# Row 93
$ws.Range("H93").Value = 1450
$ws.Range("I93").Value = 1450
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1450
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -202
$ws.Range("N93").Value = ""

# Row 122
$ws.Range("H122").Value = 3173.6365
$ws.Range("I122").Value = 2983
$ws.Range("K122").Value = 8949
$ws.Range("M122").Value = -6499

# Row 132
$ws.Range("H132").Value = 5661.1113
$ws.Range("I132").Value = 3490.8333
$ws.Range("J132").Value = 10001.667
$ws.Range("K132").Value = 10472.4999
$ws.Range("L132").Value = 30005.001
$ws.Range("M132").Value = -7942.499899999999
$ws.Range("N132").Value = -35065.001

# Row 137
$ws.Range("H137").Value = 49426
$ws.Range("J137").Value = 49426
$ws.Range("L137").Value = 49426
$ws.Range("N137").Value = -59626


$ws = $wb.Worksheets.Item("WVR")
# Row 9
$ws.Range("H9").Value = 1406
$ws.Range("I9").Value = 1406
$ws.Range("K9").Value = 1406
$ws.Range("M9").Value = -1266

# Row 100
$ws.Range("H100").Value = 1265.8462
$ws.Range("I100").Value = 1265.8462
$ws.Range("K100").Value = 2531.6924
$ws.Range("M100").Value = -1990.6924

# Row 107
$ws.Range("H107").Value = 99.5
$ws.Range("I107").Value = 99.5
$ws.Range("K107").Value = 298.5
$ws.Range("M107").Value = 1621.5

